$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.469.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.313.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.41%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.314.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.469"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.25%  "

$ws.Range("E11").Value = "  -3.45%  "

$ws.Range("E12").Value = "  -0.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.880.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.311.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.02%  "

$ws.Range("E17").Value = "  -2.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.480.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.72%  "

$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "374.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("E25").Value = "  -4.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.445.86"
$ws.Range("D26").Style = "Normal"

$ws.Range("E28").Value = "  -4.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.54%  "

$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("E32").Value = "  -3.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.66%  "

$ws.Range("E35").Value = "  -6.86%  "

$ws.Range("E36").Value = "  -5.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.01%  "

$ws.Range("E38").Value = "  -3.89%  "

$ws.Range("E39").Value = "  -2.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.343.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0731"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.28%  "

$ws.Range("E44").Value = "  -2.71%  "

$ws.Range("E45").Value = "  -4.24%  "

$ws.Range("E46").Value = "  -4.76%  "

$ws.Range("E47").Value = "  -4.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.375.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("E50").Value = "  -6.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.02%  "
